$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from B1 (header style) to C1, then set header value
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = 2

# Fill C2:C67 with the predicted t+3 values
$values = @(
    -0.3113347074289435, 0.2065752829286065, 0.1421566416881549, -0.07694397160759842, 0.1429155484052026, -0.4527683761063582, -0.2598008161228147, -0.411365857694458,
    0.3583351033494231, -0.239615135361616, -0.0872659943667467, -0.01163689161266203, 0.05201512159417407, -0.01353990189783741, 0.3474641149721687, 0.5981113559656567,
    0.03848846692182706, 0.3958525534261808, 0.2938656074185948, 0.5182927104842734, 0.3774294530484647, -0.06047098978932161, 4.575953720226746, 0.4916963026454615,
    0.3637634469251554, 0.3013268639216591, 1.068469446076402, 5.453715541498916, 0.9684969965365088, -0.19829785748661, 0.7771228137702961, 0.9117538597084132,
    -0.6548704237507494, 0.813604258629973, 0.7675934109401431, 0.7417926700373553, 0.7173287459599338, 0.5810806496389941, 0.7522695508355323, 0.5544369560187997,
    0.690163836811005, 0.7090153754814805, 0.6744744297863949, 0.6581781264025099, -1.263522051325703, -0.9793198897550003, -0.8733249578619043, -0.6423388394497972,
    -0.04870843203017888, -0.8651112230871479, -0.8651112230871479, -1.097167705592106, -0.1869482209713249, -0.9991769555261637, -0.8841450427105935, -0.9339108452677687,
    -1.128012501669319, -0.8479030779331405, -0.4977996238259934, 0.3714301207007228, -1.209516136962105, -0.7390840125550109, -0.9025893591691596, -0.1116670023084185,
    -0.7907523108579071, -0.7623398618387645
)

$rowCount = $values.Count
$arr = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i,0] = $values[$i]
}

$ws.Range("C2:C67").Value = $arr
